# Update "想去人数" (F) / "最低票价" (G) counts across all four sheets to
# reflect the newly scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1203
$ws.Range("F5").Value = 630
$ws.Range("F7").Value = 1735
$ws.Range("F12").Value = 277
$ws.Range("F13").Value = 1656
$ws.Range("G13").Value = 58
$ws.Range("F14").Value = 318
$ws.Range("F15").Value = 1368
$ws.Range("F16").Value = 772
$ws.Range("F17").Value = 309
$ws.Range("F18").Value = 648
$ws.Range("F19").Value = 12543
$ws.Range("F20").Value = 12587
$ws.Range("F21").Value = 929
$ws.Range("F23").Value = 8
$ws.Range("F24").Value = 287
$ws.Range("F26").Value = 452
$ws.Range("F27").Value = 1955
$ws.Range("F28").Value = 18
$ws.Range("F31").Value = 650

# --- 演出 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 47
$ws.Range("F8").Value = 123
$ws.Range("F10").Value = 62

# --- 本地生活 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 139

# --- 全部类型 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1203
$ws.Range("F6").Value = 630
$ws.Range("F7").Value = 139
$ws.Range("F9").Value = 1735
$ws.Range("F16").Value = 47
$ws.Range("F17").Value = 277
$ws.Range("F18").Value = 1656
$ws.Range("G18").Value = 58
$ws.Range("F19").Value = 318
$ws.Range("F20").Value = 1368
$ws.Range("F21").Value = 772
$ws.Range("F22").Value = 309
$ws.Range("F24").Value = 648
$ws.Range("F25").Value = 12543
$ws.Range("F26").Value = 12587
$ws.Range("F27").Value = 929
$ws.Range("F29").Value = 8
$ws.Range("F30").Value = 287
$ws.Range("F32").Value = 452
$ws.Range("F35").Value = 1955
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 123
$ws.Range("F41").Value = 650
$ws.Range("F42").Value = 62
